$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 9400
$ws.Cells.Item(62, 9).Value = 4542.857
$ws.Cells.Item(62, 10).Value = 17900
$ws.Cells.Item(62, 11).Value = 4542.857
$ws.Cells.Item(62, 12).Value = 17900
$ws.Cells.Item(62, 13).Value = -3918.857
$ws.Cells.Item(62, 14).Value = -19148
$ws.Cells.Item(65, 8).Value = 9400
$ws.Cells.Item(65, 9).Value = 4542.857
$ws.Cells.Item(65, 10).Value = 17900
$ws.Cells.Item(65, 11).Value = 22714.285
$ws.Cells.Item(65, 12).Value = 89500
$ws.Cells.Item(65, 13).Value = -19594.285
$ws.Cells.Item(65, 14).Value = -95740
$ws.Cells.Item(112, 8).Value = 4666.7734
$ws.Cells.Item(112, 10).Value = 4900.78
$ws.Cells.Item(112, 12).Value = 14702.34
$ws.Cells.Item(112, 14).Value = -16918.34
$ws.Cells.Item(118, 8).Value = 2345.6667
$ws.Cells.Item(118, 9).Value = 1848.75
$ws.Cells.Item(118, 10).Value = 4333.3335
$ws.Cells.Item(118, 11).Value = 5546.25
$ws.Cells.Item(118, 12).Value = 13000.0005
$ws.Cells.Item(118, 13).Value = -3889.25
$ws.Cells.Item(118, 14).Value = -16314.0005
$ws.Cells.Item(124, 8).Value = 64390
$ws.Cells.Item(124, 10).Value = 64390
$ws.Cells.Item(124, 12).Value = 64390
$ws.Cells.Item(124, 14).Value = -74210
$ws.Cells.Item(127, 8).Value = 761.2308
$ws.Cells.Item(127, 9).Value = 457.14285
$ws.Cells.Item(127, 10).Value = 1116
$ws.Cells.Item(127, 11).Value = 1371.42855
$ws.Cells.Item(127, 12).Value = 3348
$ws.Cells.Item(127, 13).Value = 3588.57145
$ws.Cells.Item(127, 14).Value = -13268
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 14).ClearContents()
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 4095.866
$ws.Cells.Item(138, 9).Value = 3400.8635
$ws.Cells.Item(138, 10).Value = 4350.7
$ws.Cells.Item(138, 11).Value = 10202.5905
$ws.Cells.Item(138, 12).Value = 13052.1
$ws.Cells.Item(138, 13).Value = -5062.5905
$ws.Cells.Item(138, 14).Value = -23332.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 95000
$ws.Cells.Item(24, 10).Value = 95000
$ws.Cells.Item(24, 12).Value = 95000
$ws.Cells.Item(24, 14).Value = -95748
$ws.Cells.Item(32, 8).Value = 22246078
$ws.Cells.Item(32, 9).Value = 35736520
$ws.Cells.Item(32, 10).Value = 26527.295
$ws.Cells.Item(32, 11).Value = 35736520
$ws.Cells.Item(32, 12).Value = 26527.295
$ws.Cells.Item(32, 13).Value = -35736233
$ws.Cells.Item(32, 14).Value = -27101.295
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 95000
$ws.Cells.Item(68, 10).Value = 95000
$ws.Cells.Item(68, 12).Value = 95000
$ws.Cells.Item(68, 14).Value = -96622
$ws.Cells.Item(71, 8).Value = 95000
$ws.Cells.Item(71, 10).Value = 95000
$ws.Cells.Item(71, 12).Value = 285000
$ws.Cells.Item(71, 14).Value = -293112
$ws.Cells.Item(76, 8).Value = 79900
$ws.Cells.Item(76, 10).Value = 79900
$ws.Cells.Item(76, 12).Value = 79900
$ws.Cells.Item(76, 14).Value = -80576
$ws.Cells.Item(79, 8).Value = 79900
$ws.Cells.Item(79, 10).Value = 79900
$ws.Cells.Item(79, 12).Value = 79900
$ws.Cells.Item(79, 14).Value = -82240
$ws.Cells.Item(100, 8).Value = 95000
$ws.Cells.Item(100, 10).Value = 95000
$ws.Cells.Item(100, 12).Value = 95000
$ws.Cells.Item(100, 14).Value = -97164
$ws.Cells.Item(101, 8).Value = 79602
$ws.Cells.Item(101, 10).Value = 79602
$ws.Cells.Item(101, 12).Value = 79602
$ws.Cells.Item(101, 14).Value = -86092
$ws.Cells.Item(104, 8).Value = 77806
$ws.Cells.Item(104, 10).Value = 77806
$ws.Cells.Item(104, 12).Value = 77806
$ws.Cells.Item(104, 14).Value = -84794
$ws.Cells.Item(132, 8).Value = 6420590.5
$ws.Cells.Item(132, 9).Value = 19002.4
$ws.Cells.Item(132, 10).Value = 10993153
$ws.Cells.Item(132, 11).Value = 57007.2
$ws.Cells.Item(132, 12).Value = 32979459
$ws.Cells.Item(132, 13).Value = -54477.2
$ws.Cells.Item(132, 14).Value = -32984519

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 83730.75
$ws.Cells.Item(94, 9).Value = 91260.82000000001
$ws.Cells.Item(94, 11).Value = 91260.82000000001
$ws.Cells.Item(94, 13).Value = -90809.82000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6994.171
$ws.Cells.Item(31, 9).Value = 1238.2142
$ws.Cells.Item(31, 10).Value = 8293.903
$ws.Cells.Item(31, 11).Value = 1238.2142
$ws.Cells.Item(31, 12).Value = 8293.903
$ws.Cells.Item(31, 13).Value = -943.2141999999999
$ws.Cells.Item(31, 14).Value = -8883.903
$ws.Cells.Item(34, 8).Value = 6994.171
$ws.Cells.Item(34, 9).Value = 1238.2142
$ws.Cells.Item(34, 10).Value = 8293.903
$ws.Cells.Item(34, 11).Value = 1238.2142
$ws.Cells.Item(34, 12).Value = 8293.903
$ws.Cells.Item(34, 13).Value = -1036.2142
$ws.Cells.Item(34, 14).Value = -8697.903
$ws.Cells.Item(86, 8).Value = 4110
$ws.Cells.Item(86, 9).Value = 4400
$ws.Cells.Item(86, 11).Value = 4400
$ws.Cells.Item(86, 13).Value = -3277
$ws.Cells.Item(89, 8).Value = 4110
$ws.Cells.Item(89, 9).Value = 4400
$ws.Cells.Item(89, 11).Value = 22000
$ws.Cells.Item(89, 13).Value = -16384
$ws.Cells.Item(100, 8).Value = 54390
$ws.Cells.Item(100, 10).Value = 54390
$ws.Cells.Item(100, 12).Value = 54390
$ws.Cells.Item(100, 14).Value = -56554
$ws.Cells.Item(105, 8).Value = 928.8889
$ws.Cells.Item(105, 9).Value = 672
$ws.Cells.Item(105, 10).Value = 1250
$ws.Cells.Item(105, 11).Value = 672
$ws.Cells.Item(105, 12).Value = 1250
$ws.Cells.Item(105, 13).Value = 1075
$ws.Cells.Item(105, 14).Value = -4744

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 441.55554
$ws.Cells.Item(98, 9).Value = 614.6
$ws.Cells.Item(98, 10).Value = 225.25
$ws.Cells.Item(98, 11).Value = 1843.8
$ws.Cells.Item(98, 12).Value = 675.75
$ws.Cells.Item(98, 13).Value = -345.8000000000002
$ws.Cells.Item(98, 14).Value = -3671.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4787.909
$ws.Cells.Item(70, 9).Value = 4738.7085
$ws.Cells.Item(70, 11).Value = 4738.7085
$ws.Cells.Item(70, 13).Value = -4468.7085
$ws.Cells.Item(73, 8).Value = 4787.909
$ws.Cells.Item(73, 9).Value = 4738.7085
$ws.Cells.Item(73, 11).Value = 4738.7085
$ws.Cells.Item(73, 13).Value = -3802.7085
$ws.Cells.Item(132, 8).Value = 55564950
$ws.Cells.Item(132, 9).Value = 111125416
$ws.Cells.Item(132, 10).Value = 4491.5557
$ws.Cells.Item(132, 11).Value = 333376248
$ws.Cells.Item(132, 12).Value = 13474.6671
$ws.Cells.Item(132, 13).Value = -333373718
$ws.Cells.Item(132, 14).Value = -18534.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3371.4443
$ws.Cells.Item(132, 9).Value = 2741.4375
$ws.Cells.Item(132, 10).Value = 3875.45
$ws.Cells.Item(132, 11).Value = 8224.3125
$ws.Cells.Item(132, 12).Value = 11626.35
$ws.Cells.Item(132, 13).Value = -5694.3125
$ws.Cells.Item(132, 14).Value = -16686.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(97, 8).Value = 84786
$ws.Cells.Item(97, 10).Value = 84786
$ws.Cells.Item(97, 12).Value = 84786
$ws.Cells.Item(97, 14).Value = -86768
$ws.Cells.Item(126, 8).Value = 1589.2963
$ws.Cells.Item(126, 9).Value = 1553.4736
$ws.Cells.Item(126, 10).Value = 1674.375
$ws.Cells.Item(126, 11).Value = 4660.4208
$ws.Cells.Item(126, 12).Value = 5023.125
$ws.Cells.Item(126, 13).Value = -2190.4208
$ws.Cells.Item(126, 14).Value = -9963.125
$ws.Cells.Item(132, 8).Value = 12684080
$ws.Cells.Item(132, 9).Value = 3055.2856
$ws.Cells.Item(132, 10).Value = 32410118
$ws.Cells.Item(132, 11).Value = 9165.856800000001
$ws.Cells.Item(132, 12).Value = 97230354
$ws.Cells.Item(132, 13).Value = -6635.856800000001
$ws.Cells.Item(132, 14).Value = -97235414
